$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# need to be forced to Text so they keep their literal (non-numeric) string value,
# then the temporary Text number-format is reset back to Normal so no extra
# cell formatting is left behind.

$ws.Range("D2").Value = "63.811.84"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.749.34"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.112"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -18.06%  "
$ws.Range("D13").Value = "3.234.63"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "63.757.77"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000155"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "2.755.60"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.565"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.172"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0930"
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.81%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "332.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0597"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.640"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0258"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  +1.13%  "
